$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the data-driven login test fixture
$ws.Name = "TestLogin"

# Header row
$ws.Range("A1").Value = "fullName"
$ws.Range("B1").Value = "password"

# Data rows (write column B values before moving to the next row's column A
# so the shared-string table is interned in the same order as the source)
$ws.Range("B2").Value = "password1"
$ws.Range("B3").Value = "password2"
$ws.Range("A2").Value = "tester 2"
$ws.Range("A3").Value = "tester 3"

# Move the selection to C4 to match the saved sheet view
$ws.Range("C4").Select()
